$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.891.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.00%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.280.97"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.23%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.60"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.48"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.14%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.28"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0908"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.23"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.94%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.963"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.41"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.628.48"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.288.98"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.811.26"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.48"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000104"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.14"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.43"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -7.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.41"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.30"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.29"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -8.37%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.92"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.15%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.61"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "165.88"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.47"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0888"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.91"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.73"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.75%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.53%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.07%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0351"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.97%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.61"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.15%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "BitcoinSV"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.72"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.23%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.04"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.31%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.226"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.94%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.24"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "111.99"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -8.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.03"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.28"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.34"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.555.50"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.37%  "
